$wb = $excel.ActiveWorkbook

# Mars sheet: fill in new time-report entries for 2015-03-24 and 2015-03-25
$wsMars = $wb.Worksheets.Item("Mars")
$wsMars.Range("C34").Value = 5
$wsMars.Range("D34").Value = "H-möte. Kodning av spellogik, vinst och lägga bricka."
$wsMars.Range("C35").Value = 3
$wsMars.Range("D35").Value = "Fungerande vinstkoll och optimering av spellogik."

# Move the active cell / selection on the Mars sheet to D35
[void]$wsMars.Range("D35").Select()

# Mars becomes the active (selected) tab instead of Översikt
$wsMars.Activate()
